$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Shift the task rows down by one (a new task "Location Enhancemets to
# Events" is inserted at row 2) and fill in the final values for every
# affected row directly. Writing the exact end-state avoids relying on
# Rows.Insert() (which would also needlessly break apart the shared
# TEXT() formula used for column C).
# ---------------------------------------------------------------------------

# Row 2 - new task
$ws.Range("A2").Value = "Location Enhancemets to Events"
$ws.Range("B2").Value = "Complete"
$ws.Range("D2").Value = 43617
$ws.Range("C2").Formula = "=TEXT(D2,""DD-MMM-YY"")"
$ws.Range("E2").Value = 32
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "Enhance all VDAB events to support GPS location."

# Row 3 - was row 2 (Enhanced Alert Support)
$ws.Range("A3").Value = "Enhanced Alert Support"
$ws.Range("B3").Value = "Complete"
$ws.Range("D3").Value = 43647
$ws.Range("C3").Formula = "=TEXT(D3,""DD-MMM-YY"")"
$ws.Range("E3").Value = 12
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "Develop enhanced alerting capabilities needed to support Alert locations."

# Row 4 - was row 3 (USGS Service Node)
$ws.Range("A4").Value = "USGS Service Node"
$ws.Range("B4").Value = "Complete"
$ws.Range("D4").Value = 43651
$ws.Range("C4").Formula = "=TEXT(D4,""DD-MMM-YY"")"
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "Develop node to read the data from the USGS Water Watch data repository"

# Row 5 - was row 4 (GLOS Buoy Service Node)
$ws.Range("A5").Value = "GLOS Buoy Service Node"
$ws.Range("B5").Value = "Active"
$ws.Range("D5").Value = 43678
$ws.Range("C5").Formula = "=TEXT(D5,""DD-MMM-YY"")"
$ws.Range("E5").Value = 20
$ws.Range("F5").Value = 0.2
$ws.Range("G5").Value = "Develop a nde to read the latest Buoy data from the GLOS Buoy System"

# Fix up number formats / styles for E2,F2 and F3 which otherwise would
# inherit the wrong look after being re-written above.
$ws.Range("E2").NumberFormat = "0"
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").VerticalAlignment = -4160

$ws.Range("F2").NumberFormat = "0%"
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").VerticalAlignment = -4160

$ws.Range("F3").NumberFormat = "0%"
$ws.Range("F3").HorizontalAlignment = -4108
$ws.Range("F3").VerticalAlignment = -4160

# Row 10 grows two more formatted (but empty) cells, and a brand-new blank
# row 12 appears (mirroring row 11) at the bottom of the template.
$ws.Range("E10").Formula = $ws.Range("E9").Formula()
$ws.Range("E10").Value = ""
$ws.Range("E10").NumberFormat = "0"
$ws.Range("E10").HorizontalAlignment = -4108
$ws.Range("E10").VerticalAlignment = -4160

$ws.Range("F10").NumberFormat = "0%"
$ws.Range("F10").HorizontalAlignment = -4108
$ws.Range("F10").VerticalAlignment = -4160

$ws.Range("C12").Value = $ws.Range("C11").Value()
$ws.Range("C12").Value = ""
$ws.Range("C12").NumberFormat = "General"
$ws.Range("C12").HorizontalAlignment = -4108

$ws.Range("D12").NumberFormat = "dd-mmm-yy"
$ws.Range("D12").HorizontalAlignment = -4131

# Widen column A slightly to fit the new task title.
$ws.Columns("A").ColumnWidth = 25.83

# Move the active selection to A5, matching where the user left off.
$ws.Range("A5").Select()

# Leave review comments on the two tasks the author flagged.
$ws.Range("A2").AddComment("Author:`n")
$ws.Range("A3").AddComment("Author:`n")

Write-Host "edit applied"
